# Insert a new record as row 262 (a new weekly price observation), pushing
# the existing rows 262..359 down to 263..360 and extending the used range
# to A1:R360.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(262).Insert()

$ws.Cells.Item(262, 1).Value = 10
$ws.Cells.Item(262, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(262, 3).Value = "La Araucanía"
$ws.Cells.Item(262, 4).Value = 44588
$ws.Cells.Item(262, 5).Value = 9
$ws.Cells.Item(262, 6).Value = 100112023
$ws.Cells.Item(262, 7).Value = "Brócoli"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 750
$ws.Cells.Item(262, 11).Value = 1000
$ws.Cells.Item(262, 12).Value = 1000
$ws.Cells.Item(262, 13).Value = 1000
$ws.Cells.Item(262, 14).Value = "$/unidad"
$ws.Cells.Item(262, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(262, 16).Value = 1000
$ws.Cells.Item(262, 17).Value = 1
$ws.Cells.Item(262, 18).Value = "Hortaliza"
